$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns (Nama, Satuan) between Kode Barang (A) and Stok (B)
$ws.Columns.Item("B").Insert()
$ws.Columns.Item("B").Insert()

# Header row
$ws.Range("A1").Value = "Kode Barang"
$ws.Range("B1").Value = "Nama"
$ws.Range("C1").Value = "Satuan"
$ws.Range("D1").Value = "Stok"

# Row 2 - Laptop
$ws.Range("A2").Value = "10004S"
$ws.Range("B2").Value = "Laptop 1"
$ws.Range("C2").Value = "pcs"
$ws.Range("D2").Value = 99

# Row 3 - Mouse
$ws.Range("A3").Value = "10003S"
$ws.Range("B3").Value = "Mouse 1"
$ws.Range("C3").Value = "pcs"
$ws.Range("D3").Value = 7

# Row 4 - Keyboard
$ws.Range("A4").Value = "10001S"
$ws.Range("B4").Value = "Keyboard 1"
$ws.Range("C4").Value = "pcs"
$ws.Range("D4").Value = 22

$ws.Range("D3").Select()
